# Designation Master.xlsx edit
# - Rename the header in A1 from "FullName" to "Name"
# - Re-fit the header columns (their best-fit width naturally shrinks now
#   that the text is shorter) and leave the cursor on C2, matching the
#   selection state the workbook was saved with.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"

[void]$ws.Columns.Item(1).AutoFit()
[void]$ws.Columns.Item(2).AutoFit()

[void]$ws.Range("C2").Select()
